# Remade dropdown field arrows. Made them in css
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting: D3:D6 gains the yellow highlight fill (same fill as C3:C6)
# while keeping its existing wrapText alignment. C3:C6 stays as-is.
$ws.Range("D3:D6").Interior.ColorIndex = 6

# --- Text updates ---
$ws.Range("D3").Value = "Правое поле:
Абсол-е поз-е,
Ширина 100% высоты блока,
flex"

$ws.Range("E4").Value = "Стрелка вверх"
$ws.Range("E5").Value = "Стрелка вправо градиентная"

# Re-setting a multi-line value can trigger row autofit; keep the original
# fixed row height (15) intact, as in the source file.
$ws.Rows(3).RowHeight = 15

# --- Selection moves to E5 ---
$ws.Range("E5").Select()
